$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 118, shifting the existing row 118 down to row 119.
$ws.Rows.Item(118).Insert()

# Copy the old row 118 (now at row 119) values are already there; we just need to
# populate the new row 118 with the new weekly entry.
$ws.Cells.Item(118, 1).Value = 10
$ws.Cells.Item(118, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(118, 3).Value = "La Araucanía"
$ws.Cells.Item(118, 4).Value = 45239
$ws.Cells.Item(118, 4).NumberFormat = $ws.Cells.Item(119, 4).NumberFormat
$ws.Cells.Item(118, 5).Value = 9
$ws.Cells.Item(118, 6).Value = "Fruta"
$ws.Cells.Item(118, 7).Value = 100107
$ws.Cells.Item(118, 8).Value = "Otros"
$ws.Cells.Item(118, 9).Value = 100107011
$ws.Cells.Item(118, 10).Value = "Tuna"
$ws.Cells.Item(118, 11).Value = "Sin especificar"
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 55
$ws.Cells.Item(118, 14).Value = 40000
$ws.Cells.Item(118, 15).Value = 40000
$ws.Cells.Item(118, 16).Value = 40000
$ws.Cells.Item(118, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(118, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(118, 19).Value = 2500
$ws.Cells.Item(118, 20).Value = 16
